# Update the "liste référence" sheet: English technology labels in columns A and C
# were reworded/shortened (rows 399-411). Column B (French labels) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    399 = "Wood (existing)";
    400 = "Wood-electric dual system (existing)";
    401 = "Wood-oil dual system (existing)";
    402 = "Heating oil furnace (existing)";
    403 = "Heating oil furnace (new)";
    404 = "Natural gas furnace (existing)";
    405 = "Natural gas furnace (new)";
    406 = "Electric baseboard (existing)";
    407 = "Ductless air source heat pump (existing)";
    408 = "Ductless HP with EBB backup";
    409 = "Ductless HP with oil backup";
    410 = "Ductless HP with wood backup";
    411 = "Ductless HP with natural gas backup";
}

# Row order matches the order the new labels were first introduced in the
# source workbook (mirrors the shared-string table layout in the target file).
$rowOrder = @(399, 400, 401, 403, 404, 406, 402, 405, 407, 408, 409, 410, 411)

foreach ($row in $rowOrder) {
    $text = $updates[$row]
    $ws.Cells.Item($row, 1).Value = $text
    $ws.Cells.Item($row, 3).Value = $text
}

$ws.Range("A400").Select() | Out-Null
